$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1047.375
$ws.Range("I32").Value = 993.75
$ws.Range("K32").Value = 993.75
$ws.Range("M32").Value = -667.75
$ws.Range("H38").Value = 1439.6
$ws.Range("I38").Value = 604
$ws.Range("J38").Value = 1996.6666
$ws.Range("K38").Value = 1812
$ws.Range("L38").Value = 5989.9998
$ws.Range("M38").Value = -1440
$ws.Range("N38").Value = -6733.9998
$ws.Range("H43").Value = 727.05884
$ws.Range("I43").Value = 621.55554
$ws.Range("J43").Value = 845.75
$ws.Range("K43").Value = 621.55554
$ws.Range("L43").Value = 845.75
$ws.Range("M43").Value = -552.55554
$ws.Range("N43").Value = -983.75
$ws.Range("H51").Value = 3754.4285
$ws.Range("I51").Value = 1671.5714
$ws.Range("J51").Value = 4795.857
$ws.Range("K51").Value = 1671.5714
$ws.Range("L51").Value = 4795.857
$ws.Range("M51").Value = -1187.5714
$ws.Range("N51").Value = -5763.857
$ws.Range("H92").Value = 1074.2
$ws.Range("I92").Value = 837.5454999999999
$ws.Range("J92").Value = 1725
$ws.Range("K92").Value = 837.5454999999999
$ws.Range("L92").Value = 1725
$ws.Range("M92").Value = 410.4545000000001
$ws.Range("N92").Value = -4221
$ws.Range("H111").Value = 2392
$ws.Range("I111").Value = 2151.1333
$ws.Range("J111").Value = 2753.3
$ws.Range("K111").Value = 6453.3999
$ws.Range("L111").Value = 8259.900000000001
$ws.Range("M111").Value = -3386.3999
$ws.Range("N111").Value = -14393.9
$ws.Range("H115").Value = 1882.7778
$ws.Range("I115").Value = 589
$ws.Range("J115").Value = 3500
$ws.Range("K115").Value = 1767
$ws.Range("L115").Value = 10500
$ws.Range("M115").Value = -200
$ws.Range("N115").Value = -13634
$ws.Range("H123").Value = 22627.143
$ws.Range("J123").Value = 22627.143
$ws.Range("L123").Value = 22627.143
$ws.Range("N123").Value = -32427.143
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 1508.1936
$ws.Range("I137").Value = 1261.3077
$ws.Range("K137").Value = 3783.9231
$ws.Range("M137").Value = -1233.9231
$ws.Range("H138").Value = 2224778
$ws.Range("I138").Value = 1019.2
$ws.Range("J138").Value = 5719256
$ws.Range("K138").Value = 3057.6
$ws.Range("L138").Value = 17157768
$ws.Range("M138").Value = 2082.4
$ws.Range("N138").Value = -17168048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1972.6666
$ws.Range("I2").Value = 1823.9445
$ws.Range("J2").Value = 2418.8333
$ws.Range("K2").Value = 1823.9445
$ws.Range("L2").Value = 2418.8333
$ws.Range("M2").Value = -1710.9445
$ws.Range("N2").Value = -2644.8333
$ws.Range("H63").Value = 2462.0557
$ws.Range("I63").Value = 1954.4667
$ws.Range("K63").Value = 1954.4667
$ws.Range("M63").Value = -1268.4667
$ws.Range("H66").Value = 2462.0557
$ws.Range("I66").Value = 1954.4667
$ws.Range("K66").Value = 9772.333500000001
$ws.Range("M66").Value = -6340.333500000001
$ws.Range("H74").Value = 46992.543
$ws.Range("I74").Value = 112323.78
$ws.Range("J74").Value = 7793.8
$ws.Range("K74").Value = 112323.78
$ws.Range("L74").Value = 7793.8
$ws.Range("M74").Value = -111449.78
$ws.Range("N74").Value = -9541.799999999999
$ws.Range("H77").Value = 46992.543
$ws.Range("I77").Value = 112323.78
$ws.Range("J77").Value = 7793.8
$ws.Range("K77").Value = 561618.9
$ws.Range("L77").Value = 38969
$ws.Range("M77").Value = -557250.9
$ws.Range("N77").Value = -47705
$ws.Range("H102").Value = 1562.2307
$ws.Range("I102").Value = 1288.625
$ws.Range("K102").Value = 1288.625
$ws.Range("M102").Value = 333.375
$ws.Range("H116").Value = 1972.6666
$ws.Range("I116").Value = 1823.9445
$ws.Range("J116").Value = 2418.8333
$ws.Range("K116").Value = 1823.9445
$ws.Range("L116").Value = 2418.8333
$ws.Range("M116").Value = 470.0554999999999
$ws.Range("N116").Value = -7006.8333
$ws.Range("H122").Value = 2567.75
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2567.75
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 7703.25
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -12603.25
$ws.Range("H135").Value = 30160
$ws.Range("J135").Value = 30160
$ws.Range("L135").Value = 30160
$ws.Range("N135").Value = -40300
$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 35000
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1972.6666
$ws.Range("I3").Value = 1823.9445
$ws.Range("J3").Value = 2418.8333
$ws.Range("K3").Value = 1823.9445
$ws.Range("L3").Value = 2418.8333
$ws.Range("M3").Value = -1709.9445
$ws.Range("N3").Value = -2646.8333
$ws.Range("H86").Value = 1743.7727
$ws.Range("I86").Value = 1641.1428
$ws.Range("J86").Value = 1923.375
$ws.Range("K86").Value = 1641.1428
$ws.Range("L86").Value = 1923.375
$ws.Range("M86").Value = -518.1428000000001
$ws.Range("N86").Value = -4169.375
$ws.Range("H89").Value = 1743.7727
$ws.Range("I89").Value = 1641.1428
$ws.Range("J89").Value = 1923.375
$ws.Range("K89").Value = 8205.714
$ws.Range("L89").Value = 9616.875
$ws.Range("M89").Value = -2589.714
$ws.Range("N89").Value = -20848.875
$ws.Range("H107").Value = 1702.16
$ws.Range("I107").Value = 1594.2632
$ws.Range("J107").Value = 2043.8334
$ws.Range("K107").Value = 1594.2632
$ws.Range("L107").Value = 2043.8334
$ws.Range("M107").Value = 325.7367999999999
$ws.Range("N107").Value = -5883.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 807
$ws.Range("I16").Value = 650.3333
$ws.Range("J16").Value = 963.6667
$ws.Range("K16").Value = 650.3333
$ws.Range("L16").Value = 963.6667
$ws.Range("M16").Value = -363.3333
$ws.Range("N16").Value = -1537.6667
$ws.Range("H31").Value = 13159421
$ws.Range("I31").Value = 24391278
$ws.Range("J31").Value = 2103.7144
$ws.Range("K31").Value = 24391278
$ws.Range("L31").Value = 2103.7144
$ws.Range("M31").Value = -24390983
$ws.Range("N31").Value = -2693.7144
$ws.Range("H34").Value = 13159421
$ws.Range("I34").Value = 24391278
$ws.Range("J34").Value = 2103.7144
$ws.Range("K34").Value = 24391278
$ws.Range("L34").Value = 2103.7144
$ws.Range("M34").Value = -24391076
$ws.Range("N34").Value = -2507.7144
$ws.Range("H59").Value = 14714.286
$ws.Range("J59").Value = 14714.286
$ws.Range("L59").Value = 14714.286
$ws.Range("N59").Value = -17004.286
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H113").Value = 807
$ws.Range("I113").Value = 650.3333
$ws.Range("J113").Value = 963.6667
$ws.Range("K113").Value = 650.3333
$ws.Range("L113").Value = 963.6667
$ws.Range("M113").Value = 1519.6667
$ws.Range("N113").Value = -5303.6667
$ws.Range("H122").Value = 1663.826
$ws.Range("I122").Value = 1481.4286
$ws.Range("J122").Value = 1947.5555
$ws.Range("K122").Value = 4444.2858
$ws.Range("L122").Value = 5842.666499999999
$ws.Range("M122").Value = -1994.2858
$ws.Range("N122").Value = -10742.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 977986.0600000001
$ws.Range("I113").Value = 1377912.6
$ws.Range("J113").Value = 387.77777
$ws.Range("K113").Value = 4133737.8
$ws.Range("L113").Value = 1163.33331
$ws.Range("M113").Value = -4131567.8
$ws.Range("N113").Value = -5503.33331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5109.115
$ws.Range("I70").Value = 4933.3335
$ws.Range("J70").Value = 5259.7856
$ws.Range("K70").Value = 4933.3335
$ws.Range("L70").Value = 5259.7856
$ws.Range("M70").Value = -4663.3335
$ws.Range("N70").Value = -5799.7856
$ws.Range("H73").Value = 5109.115
$ws.Range("I73").Value = 4933.3335
$ws.Range("J73").Value = 5259.7856
$ws.Range("K73").Value = 4933.3335
$ws.Range("L73").Value = 5259.7856
$ws.Range("M73").Value = -3997.3335
$ws.Range("N73").Value = -7131.7856
$ws.Range("H80").Value = 2127.0667
$ws.Range("I80").Value = 2131.2307
$ws.Range("J80").Value = 2100
$ws.Range("K80").Value = 2131.2307
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -1133.2307
$ws.Range("N80").Value = -4096
$ws.Range("H83").Value = 2127.0667
$ws.Range("I83").Value = 2131.2307
$ws.Range("J83").Value = 2100
$ws.Range("K83").Value = 10656.1535
$ws.Range("L83").Value = 10500
$ws.Range("M83").Value = -5664.1535
$ws.Range("N83").Value = -20484
$ws.Range("H102").Value = 1251.56
$ws.Range("I102").Value = 1156.5625
$ws.Range("J102").Value = 1420.4445
$ws.Range("K102").Value = 1156.5625
$ws.Range("L102").Value = 1420.4445
$ws.Range("M102").Value = 465.4375
$ws.Range("N102").Value = -4664.4445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 26528.26
$ws.Range("I93").Value = 584.75
$ws.Range("J93").Value = 64264.273
$ws.Range("K93").Value = 584.75
$ws.Range("L93").Value = 64264.273
$ws.Range("M93").Value = 663.25
$ws.Range("N93").Value = -66760.273
$ws.Range("H122").Value = 7500
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 3928.5715
$ws.Range("K122").Value = 60000
$ws.Range("L122").Value = 11785.7145
$ws.Range("M122").Value = -57550
$ws.Range("N122").Value = -16685.7145
$ws.Range("H132").Value = 6769.385
$ws.Range("I132").Value = 7556
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 22668
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -20138
$ws.Range("N132").Value = -20058.5
$ws.Range("H136").Value = 2492
$ws.Range("I136").Value = 1556.2858
$ws.Range("J136").Value = 3802
$ws.Range("K136").Value = 4668.857400000001
$ws.Range("L136").Value = 11406
$ws.Range("M136").Value = -2118.857400000001
$ws.Range("N136").Value = -16506

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 66930.19
$ws.Range("I122").Value = 15934.667
$ws.Range("J122").Value = 134924.22
$ws.Range("K122").Value = 47804.001
$ws.Range("L122").Value = 404772.66
$ws.Range("M122").Value = -45354.001
$ws.Range("N122").Value = -409672.66
$ws.Range("H138").Value = 34000
$ws.Range("J138").Value = 34000
$ws.Range("L138").Value = 34000
$ws.Range("N138").Value = -44280
